# Updates the crypto price-tracker sheet with the latest scraped values.
# For each affected row, column D holds the "Price" text and column E the
# "Volume(1h)" percentage text (kept as literal strings, matching the sheet's
# existing inlineStr formatting -- including the padding spaces around "%").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '27.128.88'
$ws.Range("E2").Value = '  +3.47%  '

# Row 3: Ethereum
$ws.Range("D3").Value = '1.659.21'
$ws.Range("E3").Value = '  +3.90%  '

# Row 4: TetherUSD
$ws.Range("E4").Value = '  -0.10%  '

# Row 5: BNB
$ws.Range("D5").Value = '''215.59'
$ws.Range("E5").Value = '  +1.51%  '

# Row 6: XRP
$ws.Range("E6").Value = '  +1.06%  '

# Row 7: USDC
$ws.Range("E7").Value = '  -0.11%  '

# Row 8: Cardano
$ws.Range("E8").Value = '  +2.34%  '

# Row 9: Dogecoin
$ws.Range("E9").Value = '  +1.45%  '

# Row 10: Solana
$ws.Range("D10").Value = '''19.57'
$ws.Range("E10").Value = '  +3.42%  '

# Row 11: TRON
$ws.Range("E11").Value = '  +0.85%  '

# Row 12: WrappedliquidstakedEther2.0
$ws.Range("E12").Value = '  +3.88%  '

# Row 13: WrappedEther
$ws.Range("D13").Value = '1.663.63'
$ws.Range("E13").Value = '  +4.45%  '

# Row 14: Polkadot
$ws.Range("E14").Value = '  +1.98%  '

# Row 15: Polygon
$ws.Range("E15").Value = '  +2.88%  '

# Row 16: Litecoin
$ws.Range("D16").Value = '''64.94'
$ws.Range("E16").Value = '  +2.06%  '

# Row 17: BitcoinCash
$ws.Range("D17").Value = '''241.06'
$ws.Range("E17").Value = '  +6.31%  '

# Row 18: WrappedBTC
$ws.Range("D18").Value = '27.099.88'
$ws.Range("E18").Value = '  +3.34%  '

# Row 19: Chainlink
$ws.Range("E19").Value = '  +4.06%  '

# Row 20: ShibaInu
$ws.Range("E20").Value = '  +1.47%  '

# Row 21: Dai
$ws.Range("E21").Value = '  -0.09%  '

# Row 22: Uniswap
$ws.Range("E22").Value = '  +5.41%  '

# Row 23: Toncoin
$ws.Range("E23").Value = '  +3.80%  '

# Row 24: Avalanche
$ws.Range("D24").Value = '''9.34'
$ws.Range("E24").Value = '  +4.70%  '

# Row 25: Monero
$ws.Range("D25").Value = '''145.75'
$ws.Range("E25").Value = '  +0.09%  '

# Row 27: Cosmos
$ws.Range("D27").Value = '''7.16'
$ws.Range("E27").Value = '  +2.96%  '

# Row 28: Stellar
$ws.Range("E28").Value = '  +1.03%  '

# Row 29: EthereumClassic
$ws.Range("D29").Value = '''15.87'
$ws.Range("E29").Value = '  +3.43%  '

# Row 30: Hedera
$ws.Range("D30").Value = '''0.0499'
$ws.Range("E30").Value = '  +1.21%  '

# Row 31: PancakeSwap
$ws.Range("E31").Value = '  +1.21%  '

# Row 32: Maker
$ws.Range("D32").Value = '1.529.84'
$ws.Range("E32").Value = '  +6.07%  '

# Row 33: Filecoin
$ws.Range("E33").Value = '  +3.13%  '

# Row 34: InternetComputer(DFINITY)
$ws.Range("E34").Value = '  +3.57%  '

# Row 35: LidoDAOToken
$ws.Range("E35").Value = '  +8.49%  '

# Row 36: HuobiToken
$ws.Range("E36").Value = '  -0.15%  '

# Row 37: ImmutableX
$ws.Range("E37").Value = '  +1.90%  '

# Row 38: ARBITRUM
$ws.Range("D38").Value = '''0.895'
$ws.Range("E38").Value = '  +9.34%  '

# Row 39: VeChain
$ws.Range("E39").Value = '  +2.76%  '

# Row 40: FraxShare
$ws.Range("E40").Value = '  +3.54%  '

# Row 41: PaxDollar
$ws.Range("E41").Value = '  -0.08%  '

# Row 42: MXToken
$ws.Range("E42").Value = '  +5.10%  '

# Row 43: Aave
$ws.Range("D43").Value = '''66.33'
$ws.Range("E43").Value = '  +9.80%  '

# Row 44: RocketPoolETH
$ws.Range("D44").Value = '1.798.95'
$ws.Range("E44").Value = '  +3.63%  '

# Row 45: TrustWalletToken
$ws.Range("D45").Value = '''0.773'
$ws.Range("E45").Value = '  +2.30%  '

# Row 47: Quant
$ws.Range("D47").Value = '''90.52'
$ws.Range("E47").Value = '  +3.33%  '

# Row 48: RenderToken
$ws.Range("E48").Value = '  +4.01%  '

# Row 49: BabyDogeCoin
$ws.Range("E49").Value = '  -0.39%  '

# Row 50: Algorand
$ws.Range("E50").Value = '  +3.25%  '

# Row 51: Cronos
$ws.Range("E51").Value = '  +0.64%  '
